# Updates cryptos list prices (col D) and 1h volume percentages (col E)
# on the active worksheet to reflect the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value.
# These values must stay plain text (e.g. "213.09", "0.0513", "26.641.97")
# exactly as scraped, not be reinterpreted as numbers, so each cell is
# forced to Text format before assignment and reset to the Normal style
# afterwards (the source cells carry no explicit style).
$updates = [ordered]@{
    "D2"  = "26.641.97"
    "E2"  = "  +1.20%  "
    "D3"  = "1.636.41"
    "E3"  = "  +0.87%  "
    "E4"  = "  +0.03%  "
    "D5"  = "213.09"
    "E5"  = "  +0.28%  "
    "E6"  = "  +2.55%  "
    "E7"  = "  +0.06%  "
    "E8"  = "  +2.47%  "
    "E9"  = "  +1.46%  "
    "D10" = "19.18"
    "E10" = "  +1.43%  "
    "E11" = "  +3.12%  "
    "D12" = "1.865.36"
    "E12" = "  +0.94%  "
    "D13" = "1.630.92"
    "E13" = "  +0.44%  "
    "E14" = "  +2.14%  "
    "D15" = "0.527"
    "E15" = "  +1.59%  "
    "D16" = "26.679.41"
    "E16" = "  +1.31%  "
    "D17" = "63.35"
    "E17" = "  +1.37%  "
    "E18" = "  +2.23%  "
    "D19" = "218.93"
    "E19" = "  +7.98%  "
    "E20" = "  +0.06%  "
    "E21" = "  +0.25%  "
    "D22" = "9.50"
    "E22" = "  +1.49%  "
    "D23" = "6.21"
    "E23" = "  +2.70%  "
    "E24" = "  +0.34%  "
    "D25" = "148.92"
    "E25" = "  +3.90%  "
    "E26" = "  +0.07%  "
    "E27" = "  +0.04%  "
    "E28" = "  +4.32%  "
    "D29" = "15.43"
    "E29" = "  +1.52%  "
    "D30" = "0.0513"
    "E30" = "  -2.84%  "
    "E31" = "  -0.20%  "
    "E32" = "  +4.09%  "
    "E33" = "  -0.42%  "
    "E34" = "  +0.72%  "
    "E35" = "  -1.55%  "
    "D36" = "1.197.39"
    "E36" = "  +1.38%  "
    "E37" = "  +5.67%  "
    "D38" = "0.808"
    "E38" = "  -0.05%  "
    "E39" = "  +0.03%  "
    "D40" = "0.505"
    "E40" = "  +1.99%  "
    "E41" = "  -0.48%  "
    "D42" = "5.40"
    "E42" = "  +1.15%  "
    "E43" = "  +0.56%  "
    "D44" = "1.773.37"
    "E44" = "  +0.80%  "
    "D45" = "92.18"
    "E45" = "  -1.58%  "
    "E46" = "  +1.38%  "
    "D47" = "54.77"
    "E47" = "  +1.43%  "
    "E48" = "  +0.81%  "
    "D49" = "7.64"
    "E49" = "  +4.76%  "
    "E51" = "  +0.12%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text format so values like "213.09" or "0.0513" are not
    # reinterpreted as numbers (which would drop trailing zeros / use
    # scientific notation) and keep the original leading zero look.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    # Restore default styling so no stray style index is introduced,
    # matching the unstyled D/E data cells in the source sheet.
    $cell.Style = "Normal"
}
